$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resultados")

$data = @(
    @(2, 15, 20, 26, 28, 36, 10, 10),
    @(4, 33, 35, 37, 38, 40, 8, 8)
)

$startRow = 80
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
